# Apply "full of errors extra attempts" edit.
#
# Setup sheet: update existing score values and add a second score row.
# Attempts sheet: add a second "attempt" column.
# Scores sheet: add a second "score" column and update existing values.

$wb = $excel.ActiveWorkbook

# --- Setup sheet --------------------------------------------------------
$setup = $wb.Worksheets.Item("Setup")
$setup.Range("H2").Value = 100
$setup.Range("I2").Value = 90
$setup.Range("J2").Value = 80

$setup.Range("G3").Value = 2
$setup.Range("H3").Value = 200
$setup.Range("I3").Value = 190
$setup.Range("J3").Value = 180

# --- Attempts sheet ------------------------------------------------------
# Add a new column C ("2") holding the second attempt count for each climber.
$attempts = $wb.Worksheets.Item("Attempts")
$attempts.Range("C1").Value = 2
$attempts.Range("C2").Value = 0
$attempts.Range("C3").Value = 0

# --- Scores sheet ---------------------------------------------------------
# Add a new column D ("2") holding the second score for each climber, and
# update the existing score values.
$scores = $wb.Worksheets.Item("Scores")
$scores.Range("D1").Value = 2

$scores.Range("B2").Value = 80
$scores.Range("C2").Value = 80

$scores.Range("B3").Value = 90
$scores.Range("C3").Value = 90
